# Updated cryptos list values (price + volume%) per the commit diff.
# Column D ("Price") cells store numeric-looking strings (e.g. "1.00",
# "73.277.66") as TEXT in the workbook, so each is written with a
# leading apostrophe -- the standard Excel "force text" entry method --
# to stop them being auto-coerced into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''73.277.66'
$ws.Range("E2").Value = '  +0.52%  '

$ws.Range("D3").Value = '''3.969.00'
$ws.Range("E3").Value = '  -1.91%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''609.68'
$ws.Range("E5").Value = '  +8.56%  '

$ws.Range("D6").Value = '''168.42'
$ws.Range("E6").Value = '  +11.43%  '

$ws.Range("E7").Value = '  -2.29%  '

$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("D9").Value = '''0.769'
$ws.Range("E9").Value = '  +0.47%  '

$ws.Range("D10").Value = '''0.185'
$ws.Range("E10").Value = '  +7.34%  '

$ws.Range("D11").Value = '''55.95'
$ws.Range("E11").Value = '  +3.88%  '

$ws.Range("E12").Value = '  +2.17%  '

$ws.Range("D13").Value = '''11.25'
$ws.Range("E13").Value = '  +2.53%  '

$ws.Range("D14").Value = '''4.604.64'
$ws.Range("E14").Value = '  -1.95%  '

$ws.Range("D15").Value = '''3.977.22'
$ws.Range("E15").Value = '  -1.98%  '

$ws.Range("D16").Value = '''14.17'
$ws.Range("E16").Value = '  -2.47%  '

$ws.Range("E17").Value = '  +1.95%  '

$ws.Range("D18").Value = '''20.42'
$ws.Range("E18").Value = '  -1.88%  '

$ws.Range("D19").Value = '''73.171.31'
$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("E20").Value = '  -1.06%  '

$ws.Range("D21").Value = '''439.13'
$ws.Range("E21").Value = '  -1.46%  '

$ws.Range("E22").Value = '  +8.99%  '

$ws.Range("D23").Value = '''95.59'
$ws.Range("E23").Value = '  -2.69%  '

$ws.Range("D24").Value = '''3.37'
$ws.Range("E24").Value = '  -5.04%  '

$ws.Range("D25").Value = '''14.19'
$ws.Range("E25").Value = '  -4.24%  '

$ws.Range("D26").Value = '''4.10'
$ws.Range("E26").Value = '  -6.21%  '

$ws.Range("D27").Value = '''11.05'
$ws.Range("E27").Value = '  -2.72%  '

$ws.Range("D28").Value = '''5.97'
$ws.Range("E28").Value = '  +0.47%  '

$ws.Range("D29").Value = '''10.47'
$ws.Range("E29").Value = '  -4.66%  '

$ws.Range("D30").Value = '''36.05'
$ws.Range("E30").Value = '  -3.44%  '

$ws.Range("D31").Value = '''7.78'
$ws.Range("E31").Value = '  -2.10%  '

$ws.Range("D32").Value = '''13.80'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("E33").Value = '  +13.88%  '

$ws.Range("E34").Value = '  -4.19%  '

$ws.Range("D35").Value = '''47.70'
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("D36").Value = '''70.36'
$ws.Range("E36").Value = '  +4.45%  '

$ws.Range("D37").Value = '''647.61'
$ws.Range("E37").Value = '  -5.98%  '

$ws.Range("D38").Value = '''0.429'
$ws.Range("E38").Value = '  -5.16%  '

$ws.Range("D39").Value = '''3.40'
$ws.Range("E39").Value = '  -0.73%  '

$ws.Range("D40").Value = '''1.00'
$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("D41").Value = '''0.145'
$ws.Range("E41").Value = '  -2.95%  '

$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = '''0.0484'
$ws.Range("E43").Value = '  -2.94%  '

$ws.Range("E44").Value = '  -6.08%  '

$ws.Range("E45").Value = '  -4.82%  '

$ws.Range("D46").Value = '''3.11'
$ws.Range("E46").Value = '  +31.16%  '

$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '''0.148'
$ws.Range("E47").Value = '  -2.74%  '

$ws.Range("B48").Value = 'FLOKI'
$ws.Range("C48").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D48").Value = '''0.000298'
$ws.Range("E48").Value = '  +6.35%  '

$ws.Range("D49").Value = '''3.42'
$ws.Range("E49").Value = '  +2.81%  '

$ws.Range("E50").Value = '  -5.69%  '

$ws.Range("D51").Value = '''2.99'
$ws.Range("E51").Value = '  -4.80%  '
